$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 228, pushing existing rows 228:243 down to 229:244
$ws.Rows.Item(228).Insert()

# Populate the newly inserted row 228 with the new weekly record
$ws.Cells.Item(228, 1).Value = 10
$ws.Cells.Item(228, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(228, 3).Value = "La Araucanía"
$ws.Cells.Item(228, 4).Value = 44585
$ws.Cells.Item(228, 5).Value = 9
$ws.Cells.Item(228, 6).Value = 100112017
$ws.Cells.Item(228, 7).Value = "Apio"
$ws.Cells.Item(228, 8).Value = "Americana (o)"
$ws.Cells.Item(228, 9).Value = "Primera"
$ws.Cells.Item(228, 10).Value = 45
$ws.Cells.Item(228, 11).Value = 11000
$ws.Cells.Item(228, 12).Value = 11000
$ws.Cells.Item(228, 13).Value = 11000
$ws.Cells.Item(228, 14).Value = "$/docena de matas"
$ws.Cells.Item(228, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(228, 16).Value = 1833
$ws.Cells.Item(228, 17).Value = 6
$ws.Cells.Item(228, 18).Value = "Hortaliza"
